# Add a centered "Questionnaire 19" default header to the document's
# (single) section, per the commit: "Added headers to all questionnaires
# containing the questionnaires number to help keep track of the
# questionnaire numbers after printing."

$d = $word.ActiveDocument
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)   # wdHeaderFooterPrimary -> <w:headerReference w:type="default".../>

# Insert the header text first. Writing the text before touching any
# other header/footer variant is what materializes just the single
# "default" header part (word/header1.xml) instead of also minting the
# even-page / first-page variants.
$hdr.Range.InsertAfter("Questionnaire 19")

# Style the header paragraph: built-in "Header" style, centered.
$p = $hdr.Range.Paragraphs.First
$p.Style = "Header"
$p.Alignment = 1   # wdAlignParagraphCenter

# Style the run text itself (Arial, 12pt == sz 24 half-points), leaving
# the paragraph mark's run properties untouched.
$textRange = $hdr.Range.Duplicate
[void]$textRange.MoveEnd(1, -1)
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12

Write-Output "header added"
